# Regenerate save_data "K" column (Elo K-factor) values in place of the
# previous "Strike#"-derived values. Only column G (header "K") changes;
# all other columns (TB, PC, dS0, dSF, IP, I0, IF, date) are left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row number -> new K value, taken from the recomputed s_vals.
$kValues = @{
    2  = 1
    4  = 1
    5  = 1
    6  = 2
    7  = 3
    8  = 1
    9  = 2
    10 = 2
    11 = 0
    12 = 1
    13 = 1
    14 = 2
    15 = 2
    16 = 3
    17 = 1
    18 = 1
    19 = 1
    20 = 1
    21 = 1
    22 = 2
    23 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
